$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update SCENARIO_DESC (column D) cells for the "Setup Group Emiten" menu:
# Row 2 = Tambah, Row 3 = View, Row 4 = Ubah, Row 5 = Hapus
$ws.Range("D2").Value = "Tambah Setup Group Emiten"
$ws.Range("D3").Value = "View Setup Group Emiten"
$ws.Range("D4").Value = "Ubah Setup Group Emiten"
$ws.Range("D5").Value = "Hapus Setup Group Emiten"

# Row heights shrink now that the cell text is shorter (single line instead
# of a multi-step numbered list), matching Excel's wrap-text autofit result
$ws.Rows("2").RowHeight = 30
$ws.Rows("3").RowHeight = 30
$ws.Rows("4").EntireRow.AutoFit()
$ws.Rows("5").RowHeight = 30

# Update the active selection to match the new focus cell
$ws.Range("D5").Select()
